# Weekly update: insert a new weekly price record at the top of the
# "Macroferia Regional de Talca - Brócoli" data block (row 550), pushing
# the existing rows 550-669 down to 551-670.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 550:669 down by one row to make room for the new record.
$ws.Rows.Item(550).Insert()

# Populate the newly inserted row with the latest week's data.
$ws.Range("A550").Value = 5
$ws.Range("B550").Value = "Macroferia Regional de Talca"
$ws.Range("C550").Value = "Maule"
$ws.Range("D550").Value = 45258
$ws.Range("E550").Value = 7
$ws.Range("F550").Value = 100112023
$ws.Range("G550").Value = "Brócoli"
$ws.Range("H550").Value = "Sin especificar"
$ws.Range("I550").Value = "Primera"
$ws.Range("J550").Value = 4000
$ws.Range("K550").Value = 1200
$ws.Range("L550").Value = 1200
$ws.Range("M550").Value = 1200
$ws.Range("N550").Value = "$/unidad"
$ws.Range("O550").Value = "Región del Maule"
$ws.Range("P550").Value = 1200
$ws.Range("Q550").Value = 1
$ws.Range("R550").Value = "Hortaliza"
